$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row of trade data (row 6)
$ws.Range("A6").Value = 42647.680648148147
$ws.Range("A6").NumberFormat = "m/d/yy h:mm"

$ws.Range("B6").Value = $true
$ws.Range("C6").Value = 9941.89
$ws.Range("D6").Value = 9766.58
$ws.Range("E6").Value = 18.12
$ws.Range("F6").Value = 18.77

$ws.Range("G6").Value = $false
$ws.Range("G6").NumberFormat = "m/d/yy h:mm"

$ws.Range("H6").Value = 3.59
$ws.Range("I6").Value = $false
